$d = $word.ActiveDocument

# --- Edit paragraph 2: "Responsive websites and applications is a new approach..." ---
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Responsive websites and applications </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>is</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> a new approach that is used across the developing world. The steps in making a website responsive web applications allow for all users to access a website on different devices. </w:t></w:r></w:p>')

# --- Edit paragraph 3: "The benefits of a responsive website site is that allows..." ---
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The benefits of a responsive website site </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>is</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> that allows for users to view your website on multiple different devices, and ensures that you don’t lose views due to them not be able to view a website on their phone.</w:t></w:r><w:r><w:t xml:space="preserve"> Another benefit is that there are libraries out there like bootstrap that have preset sets that can make it easier for you. From a performance </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">prospective,  it</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> takes less time to create stand-alone mobile sites, easier to monitor traffic and easier to maintain the website.</w:t></w:r></w:p>')

# --- Append new content after the last paragraph ---
$lastPara = $d.Paragraphs.Last
$lastRng = $lastPara.Range
$lastRng.Collapse(0)
$lastRng.InsertParagraphAfter()
$newP1 = $d.Paragraphs.Last
$newP1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>')

$p_after1 = $d.Paragraphs.Last
$rng_after1 = $p_after1.Range
$rng_after1.Collapse(0)
$rng_after1.InsertParagraphAfter()
$newP2 = $d.Paragraphs.Last
$newP2.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>')

$p_after2 = $d.Paragraphs.Last
$rng_after2 = $p_after2.Range
$rng_after2.Collapse(0)
$rng_after2.InsertParagraphAfter()
$newP3 = $d.Paragraphs.Last
$newP3.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">As you can see both of those web technologies are a very good tool. Even though it would be nice to have animation on this website, it would not be suitable to have on the website as it would distract the user from other important information on the page. </w:t></w:r></w:p>')

$p_after3 = $d.Paragraphs.Last
$rng_after3 = $p_after3.Range
$rng_after3.Collapse(0)
$rng_after3.InsertParagraphAfter()
$newP4 = $d.Paragraphs.Last
$newP4.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Responsive Websites are the key to having more views on your website. It will allow for all to have access to the website has every 2</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:t xml:space="preserve"> click on domain in the google search is almost done on a phone. This stat tells me that responsive website are the key </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>at the moment</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> to get viewers wanting to come and visit your website several times. I would be </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>recommend</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> that we implement the library for bootstrap, has it as many different classes we can use to make the website responsive.</w:t></w:r></w:p>')

Write-Output "DONE"
Write-Output $d.Paragraphs.Count
